# repull data, push all data, mean calculation
# Update the "dSF" column (column F) with the re-pulled values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 6
    3  = 4
    4  = -2
    6  = 4
    7  = -1
    8  = -1
    10 = -5
    11 = 4
    12 = -2
    13 = -2
    14 = 2
    15 = 3
    16 = 4
    17 = -4
    19 = 3
    20 = -2
    21 = 2
    22 = 4
    24 = 1
    25 = 1
    26 = 1
    27 = -4
    28 = -4
    29 = -5
    30 = -3
    31 = 3
    32 = 1
    33 = 2
    34 = 3
    35 = 5
    36 = -1
    38 = 7
    39 = 0
    40 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
